$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 15:22"

# --- Insert "Gipuzkoa/Guipuzcoa" as a new province row right after "Sevilla" ---
# (pushes Asturias/Segovia/Leon down one row; Gipuzkoa/Guipuzcoa's old row is reused
# by the row that used to be below Caceres, i.e. everything shifts up by one)
# Row 21 = Sevilla, Row 22 = Asturias (before) ... Row 26 = Caceres
$ws.Rows.Item(22).Insert()

# Row 22 is now a new blank row for Gipuzkoa/Guipuzcoa, with updated case numbers
$ws.Cells.Item(22, 1).Value = "Gipuzkoa/Guipuzcoa"
$ws.Cells.Item(22, 2).Value = 2328
$ws.Cells.Item(22, 3).Value = 6144
$ws.Cells.Item(22, 4).Value = 4953
$ws.Cells.Item(22, 5).Value = 209

# The old duplicate "Gipuzkoa/Guipuzcoa" row (now pushed down to row 26) is removed
$ws.Rows.Item(26).Delete()

# --- Other numeric updates across the province table ---

# Galicia (row 6): Muertes 444 -> 448
$ws.Cells.Item(6, 5).Value = 448

# Bizkaia/Vizcaya (row 7): Casos totales 6667 -> 7010, Muertes 515 -> 538
$ws.Cells.Item(7, 2).Value = 7010
$ws.Cells.Item(7, 5).Value = 538

# Navarra (row 10): Casos totales 4621 -> 4656, Recuperados 3260 -> 3295
$ws.Cells.Item(10, 2).Value = 4656
$ws.Cells.Item(10, 4).Value = 3295

# Araba/Alava (row 16): Casos totales 3156 -> 3231, Muertes 308 -> 315
$ws.Cells.Item(16, 2).Value = 3231
$ws.Cells.Item(16, 5).Value = 315

# Asturias (now row 23): Casos totales 2285 -> 2298, Casos activos 636 -> 596,
# Recuperados 1469 -> 1506, Muertes 167 -> 196
$ws.Cells.Item(23, 2).Value = 2298
$ws.Cells.Item(23, 3).Value = 596
$ws.Cells.Item(23, 4).Value = 1506
$ws.Cells.Item(23, 5).Value = 196

# Segovia (now row 24): Casos activos 1031 -> 636, Recuperados 927 -> 1469,
# Muertes 290 -> 167
$ws.Cells.Item(24, 3).Value = 636
$ws.Cells.Item(24, 4).Value = 1469
$ws.Cells.Item(24, 5).Value = 167

# Leon (now row 25): Casos totales 2266 -> 2285, Casos activos 6144 -> 1031,
# Recuperados 4953 -> 927, Muertes 197 -> 290
$ws.Cells.Item(25, 2).Value = 2285
$ws.Cells.Item(25, 3).Value = 1031
$ws.Cells.Item(25, 4).Value = 927
$ws.Cells.Item(25, 5).Value = 290
